$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values (column G) computed from regenerated save_data (K instead of Strike#)
$kValues = @{
    2 = 1
    3 = 3
    4 = 2
    5 = 1
    6 = 1
    7 = 2
    8 = 3
    9 = 2
    10 = 1
    11 = 0
    12 = 1
    13 = 0
    14 = 1
    15 = 2
    16 = 1
    17 = 1
    18 = 2
    19 = 1
    20 = 4
    21 = 1
    22 = 1
    23 = 1
    24 = 2
    25 = 1
    26 = 3
    27 = 2
    28 = 2
    29 = 0
    30 = 1
    31 = 1
    32 = 3
    33 = 3
    34 = 1
    35 = 2
    36 = 0
    37 = 0
    38 = 0
    39 = 1
    40 = 0
    41 = 1
    42 = 0
    43 = 0
    44 = 2
    45 = 3
    46 = 2
    47 = 2
    48 = 1
    49 = 1
    50 = 1
    51 = 2
    52 = 1
    53 = 1
    54 = 2
    55 = 0
    56 = 0
    57 = 2
    58 = 1
    59 = 1
    60 = 2
    61 = 1
    62 = 1
    63 = 1
    64 = 3
    65 = 2
    66 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}

$wb.Save()
